$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.373.15"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "1.621.76"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'212.26"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.248"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").Value = "1.850.39"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").Value = "1.629.86"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'4.09"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  -0.80%  "
$ws.Range("D16").Value = "'64.05"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "26.393.60"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "'213.97"
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").Value = "'6.19"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").Value = "'9.26"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("D25").Value = "'147.59"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.84%  "
$ws.Range("D28").Value = "'6.81"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "'15.54"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("D32").Value = "'3.31"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("E34").Value = "  -1.23%  "
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").Value = "1.211.51"
$ws.Range("E36").Value = "  +3.88%  "
$ws.Range("E37").Value = "  +3.16%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'0.793"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("D40").Value = "'0.498"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").Value = "'2.25"
$ws.Range("E41").Value = "  -3.14%  "
$ws.Range("D42").Value = "'0.791"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").Value = "1.758.90"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "'92.55"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'54.57"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").Value = "'7.60"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").Value = "'0.408"
$ws.Range("E51").Value = "  -0.51%  "
